$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "248.24"
$ws.Cells.Item(2,7).NumberFormat = "@"
$ws.Cells.Item(2,7).Value = "8"
$ws.Cells.Item(3,7).NumberFormat = "@"
$ws.Cells.Item(3,7).Value = "8"
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "5.555"
$ws.Cells.Item(4,7).NumberFormat = "@"
$ws.Cells.Item(4,7).Value = "8"
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "0.05618"
$ws.Cells.Item(5,7).NumberFormat = "@"
$ws.Cells.Item(5,7).Value = "8"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "3.401"
$ws.Cells.Item(6,7).NumberFormat = "@"
$ws.Cells.Item(6,7).Value = "8"
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "6.486"
$ws.Cells.Item(7,7).NumberFormat = "@"
$ws.Cells.Item(7,7).Value = "8"
$ws.Cells.Item(8,2).Value = "FTXToken"
$ws.Cells.Item(8,3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "1.072"
$ws.Cells.Item(8,5).Value = "7FTXTokenFTT"
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = "8"
$ws.Cells.Item(9,2).Value = "MXToken"
$ws.Cells.Item(9,3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.8015"
$ws.Cells.Item(9,5).Value = "8MXTokenMX"
$ws.Cells.Item(9,7).NumberFormat = "@"
$ws.Cells.Item(9,7).Value = "8"
$ws.Cells.Item(10,2).Value = "One"
$ws.Cells.Item(10,3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.01172"
$ws.Cells.Item(10,5).Value = "9OneONEBestin24h"
$ws.Cells.Item(10,7).NumberFormat = "@"
$ws.Cells.Item(10,7).Value = "8"
$ws.Cells.Item(11,2).Value = "WazirX"
$ws.Cells.Item(11,3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.1427"
$ws.Cells.Item(11,5).Value = "10WazirXWRX"
$ws.Cells.Item(11,7).NumberFormat = "@"
$ws.Cells.Item(11,7).Value = "8"
$ws.Cells.Item(12,2).Value = "MandalaExchangeToken"
$ws.Cells.Item(12,3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "0.07334"
$ws.Cells.Item(12,5).Value = "11MandalaExchangeTokenMDX"
$ws.Cells.Item(12,7).NumberFormat = "@"
$ws.Cells.Item(12,7).Value = "8"
$ws.Cells.Item(13,2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "0.03194"
$ws.Cells.Item(13,5).Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Cells.Item(13,7).NumberFormat = "@"
$ws.Cells.Item(13,7).Value = "8"
$ws.Cells.Item(14,2).Value = "BitrueCoin"
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "0.02969"
$ws.Cells.Item(14,5).Value = "13BitrueCoinBTR"
$ws.Cells.Item(14,7).NumberFormat = "@"
$ws.Cells.Item(14,7).Value = "8"
$ws.Cells.Item(15,2).Value = "BitMartToken"
$ws.Cells.Item(15,3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "0.09265"
$ws.Cells.Item(15,5).Value = "14BitMartTokenBMX"
$ws.Cells.Item(15,7).NumberFormat = "@"
$ws.Cells.Item(15,7).Value = "8"
$ws.Cells.Item(16,2).Value = "BitForexToken"
$ws.Cells.Item(16,3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "0.001662"
$ws.Cells.Item(16,5).Value = "15BitForexTokenBF"
$ws.Cells.Item(16,7).NumberFormat = "@"
$ws.Cells.Item(16,7).Value = "8"
$ws.Cells.Item(17,2).Value = "CoinExToken"
$ws.Cells.Item(17,3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "0.04712"
$ws.Cells.Item(17,5).Value = "16CoinExTokenCET"
$ws.Cells.Item(17,7).NumberFormat = "@"
$ws.Cells.Item(17,7).Value = "8"
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "0.006264"
$ws.Cells.Item(18,7).NumberFormat = "@"
$ws.Cells.Item(18,7).Value = "8"
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "0.001053"
$ws.Cells.Item(19,7).NumberFormat = "@"
$ws.Cells.Item(19,7).Value = "8"
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "0.003831"
$ws.Cells.Item(20,7).NumberFormat = "@"
$ws.Cells.Item(20,7).Value = "8"
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "0.0001500"
$ws.Cells.Item(21,7).NumberFormat = "@"
$ws.Cells.Item(21,7).Value = "8"
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "0.0004601"
$ws.Cells.Item(22,7).NumberFormat = "@"
$ws.Cells.Item(22,7).Value = "8"
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "3.983"
$ws.Cells.Item(23,7).NumberFormat = "@"
$ws.Cells.Item(23,7).Value = "8"
$ws.Cells.Item(24,7).NumberFormat = "@"
$ws.Cells.Item(24,7).Value = "8"
$ws.Cells.Item(25,7).NumberFormat = "@"
$ws.Cells.Item(25,7).Value = "8"
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "0.1277"
$ws.Cells.Item(26,7).NumberFormat = "@"
$ws.Cells.Item(26,7).Value = "8"
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "2.588"
$ws.Cells.Item(27,7).NumberFormat = "@"
$ws.Cells.Item(27,7).Value = "8"
$ws.Cells.Item(28,7).NumberFormat = "@"
$ws.Cells.Item(28,7).Value = "8"
$ws.Cells.Item(29,7).NumberFormat = "@"
$ws.Cells.Item(29,7).Value = "8"
$ws.Cells.Item(30,7).NumberFormat = "@"
$ws.Cells.Item(30,7).Value = "8"
$ws.Cells.Item(31,7).NumberFormat = "@"
$ws.Cells.Item(31,7).Value = "8"
$ws.Cells.Item(32,7).NumberFormat = "@"
$ws.Cells.Item(32,7).Value = "8"
$ws.Cells.Item(33,7).NumberFormat = "@"
$ws.Cells.Item(33,7).Value = "8"
$ws.Cells.Item(34,7).NumberFormat = "@"
$ws.Cells.Item(34,7).Value = "8"
$ws.Cells.Item(35,7).NumberFormat = "@"
$ws.Cells.Item(35,7).Value = "8"
$ws.Cells.Item(36,7).NumberFormat = "@"
$ws.Cells.Item(36,7).Value = "8"
$ws.Cells.Item(37,7).NumberFormat = "@"
$ws.Cells.Item(37,7).Value = "8"
$ws.Cells.Item(38,7).NumberFormat = "@"
$ws.Cells.Item(38,7).Value = "8"
$ws.Cells.Item(39,7).NumberFormat = "@"
$ws.Cells.Item(39,7).Value = "8"
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "0.04189"
$ws.Cells.Item(40,7).NumberFormat = "@"
$ws.Cells.Item(40,7).Value = "8"
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "0.007002"
$ws.Cells.Item(41,7).NumberFormat = "@"
$ws.Cells.Item(41,7).Value = "8"
$ws.Cells.Item(42,2).Value = "BKEXToken"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.1045"
$ws.Cells.Item(42,5).Value = "41BKEXTokenBKK"
$ws.Cells.Item(42,7).NumberFormat = "@"
$ws.Cells.Item(42,7).Value = "8"
$ws.Cells.Item(43,2).Value = "CEJI"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "0.003301"
$ws.Cells.Item(43,5).Value = "42CEJICEJI"
$ws.Cells.Item(43,7).NumberFormat = "@"
$ws.Cells.Item(43,7).Value = "8"
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "0.008698"
$ws.Cells.Item(44,7).NumberFormat = "@"
$ws.Cells.Item(44,7).Value = "8"
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "0.00005636"
$ws.Cells.Item(45,7).NumberFormat = "@"
$ws.Cells.Item(45,7).Value = "8"
$ws.Cells.Item(46,7).NumberFormat = "@"
$ws.Cells.Item(46,7).Value = "8"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "0.6801"
$ws.Cells.Item(47,7).NumberFormat = "@"
$ws.Cells.Item(47,7).Value = "8"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "0.02763"
$ws.Cells.Item(48,7).NumberFormat = "@"
$ws.Cells.Item(48,7).Value = "8"
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "0.00002100"
$ws.Cells.Item(49,7).NumberFormat = "@"
$ws.Cells.Item(49,7).Value = "8"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "0.01010"
$ws.Cells.Item(50,7).NumberFormat = "@"
$ws.Cells.Item(50,7).Value = "8"
$ws.Cells.Item(51,7).NumberFormat = "@"
$ws.Cells.Item(51,7).Value = "8"
